# reporte_plantanova.xlsx - "cambios en la tabla de excel"
#
# On the PROCESO sheet the three mini-tables (header rows 1, 11, 23) share the
# same column headers in E:G -> Viabilidad, Siembra, Germinacion.
# The column order is changed to: Siembra, Germinacion, Viabilidad
# (i.e. the "Viabilidad" column is moved after "Germinacion").
# The header text AND the data values that sit underneath must move together.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PROCESO")

# Header rows for each of the three blocks, and the data rows belonging to
# each block.
$headerRows = @(1, 11, 23)
$dataRows = @{
    1  = @(2, 3)
    11 = @(12, 13, 14, 15)
    23 = @(24, 25, 26, 27)
}

foreach ($headerRow in $headerRows) {
    # Re-order the headers: E<-Siembra, F<-Germinacion, G<-Viabilidad
    $ws.Cells.Item($headerRow, 5).Value = "Siembra"
    $ws.Cells.Item($headerRow, 6).Value = "Germinacion"
    $ws.Cells.Item($headerRow, 7).Value = "Viabilidad"

    foreach ($r in $dataRows[$headerRow]) {
        $oldE = $ws.Cells.Item($r, 5).Value()
        $oldF = $ws.Cells.Item($r, 6).Value()
        $oldG = $ws.Cells.Item($r, 7).Value()

        # rotate left: new E = old Siembra (was F), new F = old Germinacion (was G), new G = old Viabilidad (was E)
        if ($oldF -eq $null) { $ws.Cells.Item($r, 5).ClearContents() } else { $ws.Cells.Item($r, 5).Value = $oldF }
        if ($oldG -eq $null) { $ws.Cells.Item($r, 6).ClearContents() } else { $ws.Cells.Item($r, 6).Value = $oldG }
        if ($oldE -eq $null) { $ws.Cells.Item($r, 7).ClearContents() } else { $ws.Cells.Item($r, 7).Value = $oldE }
    }
}
